# Auto-generated Excel COM-interop script to apply scheduled-runner value updates
# across multiple worksheets (ALC, ARM, BSM, CRP, GSM, LTW, WVR).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 259.75
$ws.Range("J33").Value = 0
$ws.Range("L33").Value = 0
$ws.Range("N33").ClearContents()
$ws.Range("H53").Value = 402.23077
$ws.Range("J53").Value = 433.91666
$ws.Range("L53").Value = 433.91666
$ws.Range("N53").Value = -1707.91666
$ws.Range("H88").Value = 503276.22
$ws.Range("I88").Value = 1126620.1
$ws.Range("K88").Value = 1126620.1
$ws.Range("M88").Value = -1126214.1
$ws.Range("H91").Value = 503276.22
$ws.Range("I91").Value = 1126620.1
$ws.Range("K91").Value = 1126620.1
$ws.Range("M91").Value = -1125216.1
$ws.Range("H96").Value = 417.54544
$ws.Range("I96").Value = 383.125
$ws.Range("J96").Value = 509.33334
$ws.Range("K96").Value = 1149.375
$ws.Range("L96").Value = 1528.00002
$ws.Range("M96").Value = 223.625
$ws.Range("N96").Value = -4274.000019999999
$ws.Range("H100").Value = 2133.3333
$ws.Range("I100").Value = 1314.2858
$ws.Range("K100").Value = 1314.2858
$ws.Range("M100").Value = -773.2858000000001
$ws.Range("H113").Value = 3488.0435
$ws.Range("I113").Value = 2492.2727
$ws.Range("J113").Value = 4400.8335
$ws.Range("K113").Value = 2492.2727
$ws.Range("L113").Value = 4400.8335
$ws.Range("M113").Value = 761.7273
$ws.Range("N113").Value = -10908.8335
$ws.Range("H114").Value = 92906.78
$ws.Range("J114").Value = 92906.78
$ws.Range("L114").Value = 92906.78
$ws.Range("N114").Value = -101584.78

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 267.63635
$ws.Range("I4").Value = 43.75
$ws.Range("K4").Value = 43.75
$ws.Range("M4").Value = 72.25
$ws.Range("H32").Value = 6783.4033
$ws.Range("I32").Value = 2121.244
$ws.Range("K32").Value = 2121.244
$ws.Range("M32").Value = -1834.244
$ws.Range("H45").Value = 10587235
$ws.Range("I45").Value = 3970
$ws.Range("K45").Value = 3970
$ws.Range("M45").Value = -3593
$ws.Range("H122").Value = 2409.8823
$ws.Range("I122").Value = 1621.9
$ws.Range("J122").Value = 3535.5715
$ws.Range("K122").Value = 4865.700000000001
$ws.Range("L122").Value = 10606.7145
$ws.Range("M122").Value = -2415.700000000001
$ws.Range("N122").Value = -15506.7145

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 204473.31
$ws.Range("I20").Value = 204473.31
$ws.Range("J20").Value = 0
$ws.Range("K20").Value = 204473.31
$ws.Range("L20").Value = 0
$ws.Range("M20").Value = -204226.31
$ws.Range("N20").ClearContents()
$ws.Range("H94").Value = 5084.5
$ws.Range("I94").Value = 3556
$ws.Range("J94").Value = 9670
$ws.Range("K94").Value = 3556
$ws.Range("L94").Value = 9670
$ws.Range("M94").Value = -3105
$ws.Range("N94").Value = -10572

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 223.20833
$ws.Range("I7").Value = 95.42856999999999
$ws.Range("K7").Value = 95.42856999999999
$ws.Range("M7").Value = 17.57143000000001
$ws.Range("H31").Value = 2644.975
$ws.Range("I31").Value = 2066.0356
$ws.Range("J31").Value = 3995.8333
$ws.Range("K31").Value = 2066.0356
$ws.Range("L31").Value = 3995.8333
$ws.Range("M31").Value = -1771.0356
$ws.Range("N31").Value = -4585.8333
$ws.Range("H34").Value = 2644.975
$ws.Range("I34").Value = 2066.0356
$ws.Range("J34").Value = 3995.8333
$ws.Range("K34").Value = 2066.0356
$ws.Range("L34").Value = 3995.8333
$ws.Range("M34").Value = -1864.0356
$ws.Range("N34").Value = -4399.8333
$ws.Range("H62").Value = 5340.25
$ws.Range("I62").Value = 5864
$ws.Range("J62").Value = 3769
$ws.Range("K62").Value = 5864
$ws.Range("L62").Value = 3769
$ws.Range("M62").Value = -5240
$ws.Range("N62").Value = -5017
$ws.Range("H65").Value = 5340.25
$ws.Range("I65").Value = 5864
$ws.Range("J65").Value = 3769
$ws.Range("K65").Value = 29320
$ws.Range("L65").Value = 18845
$ws.Range("M65").Value = -26200
$ws.Range("N65").Value = -25085
$ws.Range("H134").Value = 2499859
$ws.Range("I134").Value = 3107844.2
$ws.Range("K134").Value = 9323532.600000001
$ws.Range("M134").Value = -9320997.600000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 342.57144
$ws.Range("I2").Value = 342.8889
$ws.Range("J2").Value = 342
$ws.Range("K2").Value = 342.8889
$ws.Range("L2").Value = 342
$ws.Range("M2").Value = -229.8889
$ws.Range("N2").Value = -568
$ws.Range("H70").Value = 5556.5
$ws.Range("I70").Value = 5569.143
$ws.Range("J70").Value = 5538.8
$ws.Range("K70").Value = 5569.143
$ws.Range("L70").Value = 5538.8
$ws.Range("M70").Value = -5299.143
$ws.Range("N70").Value = -6078.8
$ws.Range("H73").Value = 5556.5
$ws.Range("I73").Value = 5569.143
$ws.Range("J73").Value = 5538.8
$ws.Range("K73").Value = 5569.143
$ws.Range("L73").Value = 5538.8
$ws.Range("M73").Value = -4633.143
$ws.Range("N73").Value = -7410.8
$ws.Range("H102").Value = 1208.3334
$ws.Range("J102").Value = 1249.5
$ws.Range("L102").Value = 1249.5
$ws.Range("N102").Value = -4493.5
$ws.Range("H109").Value = 51995.453
$ws.Range("J109").Value = 51995.453
$ws.Range("L109").Value = 51995.453
$ws.Range("N109").Value = -54075.453
$ws.Range("H122").Value = 9171108
$ws.Range("I122").Value = 11004580
$ws.Range("K122").Value = 33013740
$ws.Range("M122").Value = -33011290
$ws.Range("H125").Value = 50326
$ws.Range("J125").Value = 50326
$ws.Range("L125").Value = 50326
$ws.Range("N125").Value = -55246
$ws.Range("H132").Value = 3565.5483
$ws.Range("I132").Value = 3204.55
$ws.Range("J132").Value = 4221.909
$ws.Range("K132").Value = 9613.650000000001
$ws.Range("L132").Value = 12665.727
$ws.Range("M132").Value = -7083.650000000001
$ws.Range("N132").Value = -17725.727
$ws.Range("H141").Value = 147333.33
$ws.Range("I141").Value = 72000
$ws.Range("K141").Value = 72000
$ws.Range("M141").Value = -66820

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1251.6316
$ws.Range("I22").Value = 1434.7273
$ws.Range("J22").Value = 999.875
$ws.Range("K22").Value = 1434.7273
$ws.Range("L22").Value = 999.875
$ws.Range("M22").Value = -1139.7273
$ws.Range("N22").Value = -1589.875
$ws.Range("H27").Value = 1251.6316
$ws.Range("I27").Value = 1434.7273
$ws.Range("J27").Value = 999.875
$ws.Range("K27").Value = 1434.7273
$ws.Range("L27").Value = 999.875
$ws.Range("M27").Value = -1327.7273
$ws.Range("N27").Value = -1213.875
$ws.Range("H40").Value = 11115754
$ws.Range("I40").Value = 5000
$ws.Range("K40").Value = 5000
$ws.Range("M40").Value = -4864
$ws.Range("H55").Value = 6061212
$ws.Range("I55").Value = 173.46666
$ws.Range("J55").Value = 11112077
$ws.Range("K55").Value = 173.46666
$ws.Range("L55").Value = 11112077
$ws.Range("M55").Value = -0.4666599999999903
$ws.Range("N55").Value = -11112423
$ws.Range("H82").Value = 2298
$ws.Range("I82").Value = 2298
$ws.Range("J82").Value = 0
$ws.Range("K82").Value = 2298
$ws.Range("L82").Value = 0
$ws.Range("M82").Value = -1937
$ws.Range("N82").ClearContents()
$ws.Range("H85").Value = 2298
$ws.Range("I85").Value = 2298
$ws.Range("J85").Value = 0
$ws.Range("K85").Value = 2298
$ws.Range("L85").Value = 0
$ws.Range("M85").Value = -1050
$ws.Range("N85").ClearContents()
$ws.Range("H122").Value = 81253350
$ws.Range("I122").Value = 83336780
$ws.Range("K122").Value = 250010340
$ws.Range("M122").Value = -250007890

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2188.375
$ws.Range("I122").Value = 1514.375
$ws.Range("J122").Value = 2862.375
$ws.Range("K122").Value = 4543.125
$ws.Range("L122").Value = 8587.125
$ws.Range("M122").Value = -2093.125
$ws.Range("N122").Value = -13487.125
$ws.Range("H132").Value = 1880.0625
$ws.Range("I132").Value = 1234.7273
$ws.Range("J132").Value = 3299.8
$ws.Range("K132").Value = 3704.1819
$ws.Range("L132").Value = 9899.400000000001
$ws.Range("M132").Value = -1174.1819
$ws.Range("N132").Value = -14959.4

